# Commit: "added human readable metadata"
#
# The "notes" column (L) on Sheet1 is merged with the adjacent, now-removed
# "M" column: wherever column M held a human-readable note, that note
# becomes the value of column L (overwriting whatever was there before —
# a handful of rows had stray "actual start date: ..." placeholders in L
# that get replaced by the real note), and column M is then deleted
# entirely since it's redundant.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Sheet1 data spans rows 1-104; column M is column 13, column L is column 12.
For ($r = 1; $r -le 104; $r++) {
    $mVal = $ws.Cells.Item($r, 13).Value2
    if ($mVal -ne $null) {
        $ws.Cells.Item($r, 12).Value = $mVal
    }
}

# Column M is now fully absorbed into column L - remove it.
$ws.Columns.Item(13).Delete()

# Column K ("More Info" links) picks up a best-fit width once M is gone.
$ws.Columns.Item(11).ColumnWidth = 17.5

# Reflect where the user ended up looking: selecting the newly-consolidated
# notes column and scrolled down a bit.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("L29:L87").Select()
